# V1.2 - graphs excel
# Update the README "how to create an app password" instructions (step 5)
# with the extra note about removing spaces, and restore the sheet
# selections left over from editing (README -> B16, emails -> G14), while
# keeping "emails" as the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- README sheet -----------------------------------------------------
$wsReadme = $wb.Worksheets.Item("README")
$wsReadme.Range("B15").Value = "5) Copiar e colar na célula da planilha. Caso possua espaços, é necessário retirar"

# Leave the cursor on B16 on this sheet (matches the saved view state),
# then hop back to "emails" so it remains the active tab.
$wsReadme.Activate()
$wsReadme.Range("B16").Select()

# --- emails sheet -------------------------------------------------------
$wsEmails = $wb.Worksheets.Item("emails")
$wsEmails.Activate()
$wsEmails.Range("G14").Select()
